$wb = $excel.ActiveWorkbook

# Overview sheet: Status -> "Ready for handoff", Latest Handoff Date updated
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-25-11 20:25:33"

# zh-cn sheet: Status -> "Ready for handoff", Latest Handoff Datetime updated
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("E2").Value = "2016-03-11 20:25:30"

# de-de sheet: Status -> "Ready for handoff", Latest Handoff Datetime updated
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("E2").Value = "2016-03-11 20:25:33"
